# This script inserts 4 new weekly price records (rows 645-648) into the
# "Fruta, Terminal Hortofrutícola Agro Chillán - Naranja" sheet, pushing
# the existing rows 645-725 down to 649-729 (matching the commit
# "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the current row 645, shifting the
# rest of the table (old rows 645-725) down to rows 649-729.
$ws.Range("A645:A648").EntireRow.Insert()

# Helper values that are constant across every data row in this sheet.
$marketId   = 7
$market     = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$productId  = 100102
$producto   = "Cítricos"
$catId      = 100102005
$categoria  = "Naranja"
$unidad     = "$/bandeja 15 kilos granel"
$origen     = "Región de O'Higgins"
$kgUnidad   = 15

# New row 645: Fukumoto / Especial
$r = 645
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45131
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Fukumoto"
$ws.Cells.Item($r,12).Value = "Especial"
$ws.Cells.Item($r,13).Value = 60
$ws.Cells.Item($r,14).Value = 10000
$ws.Cells.Item($r,15).Value = 10000
$ws.Cells.Item($r,16).Value = 10000
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 667
$ws.Cells.Item($r,20).Value = $kgUnidad

# New row 646: Fukumoto / Primera
$r = 646
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45131
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Fukumoto"
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 100
$ws.Cells.Item($r,14).Value = 8000
$ws.Cells.Item($r,15).Value = 8000
$ws.Cells.Item($r,16).Value = 8000
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 533
$ws.Cells.Item($r,20).Value = $kgUnidad

# New row 647: Lane Late / Primera
$r = 647
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45131
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Lane Late"
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 60
$ws.Cells.Item($r,14).Value = 7000
$ws.Cells.Item($r,15).Value = 7000
$ws.Cells.Item($r,16).Value = 7000
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 467
$ws.Cells.Item($r,20).Value = $kgUnidad

# New row 648: Lane Late / Segunda
$r = 648
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45131
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Lane Late"
$ws.Cells.Item($r,12).Value = "Segunda"
$ws.Cells.Item($r,13).Value = 80
$ws.Cells.Item($r,14).Value = 6000
$ws.Cells.Item($r,15).Value = 6000
$ws.Cells.Item($r,16).Value = 6000
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 400
$ws.Cells.Item($r,20).Value = $kgUnidad
